$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.140776753425598
$ws.Range("B1").Value = 2.566413640975952
$ws.Range("C1").Value = 9.656722068786621
$ws.Range("D1").Value = 2.145643949508667
$ws.Range("E1").Value = 1.254300236701965
